# Update gh-pages generated output data (scraped attendance / price numbers
# and one cover-image URL refresh) for the 杭州-漫展信息 workbook.

$wb = $excel.ActiveWorkbook

function Set-CellValue($SheetName, $CellRef, $Value) {
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Range($CellRef).Value = $Value
}

# ---- Sheet: 展览 ----
Set-CellValue "展览" "F3"  2560
Set-CellValue "展览" "F5"  913
Set-CellValue "展览" "F7"  1356
Set-CellValue "展览" "F8"  1716
Set-CellValue "展览" "F9"  188
Set-CellValue "展览" "F11" 2392
Set-CellValue "展览" "F12" 504
Set-CellValue "展览" "F13" 169
Set-CellValue "展览" "F14" 52
Set-CellValue "展览" "F16" 110
Set-CellValue "展览" "F18" 8694
Set-CellValue "展览" "F20" 6783
Set-CellValue "展览" "F21" 10957
Set-CellValue "展览" "F24" 216
Set-CellValue "展览" "F25" 299
Set-CellValue "展览" "F26" 529
Set-CellValue "展览" "F27" 2415
Set-CellValue "展览" "F29" 179
Set-CellValue "展览" "F30" 2234
Set-CellValue "展览" "F31" 141
Set-CellValue "展览" "F33" 4453
Set-CellValue "展览" "F34" 498
Set-CellValue "展览" "F35" 7
Set-CellValue "展览" "F36" 447

# ---- Sheet: 演出 ----
Set-CellValue "演出" "F2"  8
Set-CellValue "演出" "F8"  1179
Set-CellValue "演出" "G13" 280
Set-CellValue "演出" "F15" 2

# ---- Sheet: 本地生活 ----
Set-CellValue "本地生活" "F3" 615
Set-CellValue "本地生活" "I3" "//i2.hdslb.com/bfs/openplatform/202410/d6ORQLlE1728627123348.jpeg"
Set-CellValue "本地生活" "F5" 61

# ---- Sheet: 全部类型 ----
Set-CellValue "全部类型" "F3"  615
Set-CellValue "全部类型" "I3"  "//i2.hdslb.com/bfs/openplatform/202410/d6ORQLlE1728627123348.jpeg"
Set-CellValue "全部类型" "F5"  61
Set-CellValue "全部类型" "F6"  8
Set-CellValue "全部类型" "F7"  2560
Set-CellValue "全部类型" "F9"  913
Set-CellValue "全部类型" "F11" 1356
Set-CellValue "全部类型" "F13" 1716
Set-CellValue "全部类型" "F15" 188
Set-CellValue "全部类型" "F17" 504
Set-CellValue "全部类型" "F18" 169
Set-CellValue "全部类型" "F19" 52
Set-CellValue "全部类型" "F21" 110
Set-CellValue "全部类型" "F23" 8694
Set-CellValue "全部类型" "F25" 6783
Set-CellValue "全部类型" "F26" 10957
Set-CellValue "全部类型" "F30" 216
Set-CellValue "全部类型" "F31" 299
Set-CellValue "全部类型" "F33" 529
Set-CellValue "全部类型" "G35" 280
Set-CellValue "全部类型" "F38" 4453
Set-CellValue "全部类型" "F40" 2
Set-CellValue "全部类型" "F45" 447
